$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Add a new tracked file (b9ee350c-e115-4378-a1b2-248be2ee1af7.md) as a new
# row across the Overview / zh-cn / de-de report tables (handback report).
# ---------------------------------------------------------------------------

$commitSha = "d1303fd84b6a71389dbbb699ef05f45086c563e0"
$newFile = "b9ee350c-e115-4378-a1b2-248be2ee1af7.md"
$newFileUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitSha/e2e/$newFile"

$hyperlinkColor = 15570276   # RGB(100,149,237) = FF6495ED, matches existing HyperLink style
$dateFormat = "yyyy-mm-dd HH:mm:ss"

function Style-AsHyperlink($range) {
    $range.Font.Name = "Calibri"
    $range.Font.Color = $hyperlinkColor
    $range.Font.Underline = 2
}

function Style-AsDate($range) {
    $range.NumberFormat = $dateFormat
}

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = $newFile
$wsOverview.Range("B3").Value = "e2e\$newFile"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-16 08:39:53"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $newFileUrl, "", "", "e2e\$newFile") | Out-Null
Style-AsHyperlink $wsOverview.Range("B3")
Style-AsDate $wsOverview.Range("G3")

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("A3").Value = $newFile
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "False"
$wsZh.Range("G3").Value = "b9ee350c-e115-4378-a1b2-248be2ee1af7.815491ef7fccedc18cdde3a25d641bd4f82b7eda.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-08-16 08:39:47"
$wsZh.Range("I3").Value = ""
$wsZh.Range("J3").Value = ""
$wsZh.Range("K3").Value = "0001-01-01 00:00:00"
$wsZh.Range("L3").Value = ""
$wsZh.Range("M3").Value = "True"
$wsZh.Range("N3").Value = ""
$wsZh.Range("O3").Value = "False"
$wsZh.Range("P3").Value = ""

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $newFileUrl, "", "", $newFile) | Out-Null
Style-AsHyperlink $wsZh.Range("A3")
Style-AsDate $wsZh.Range("H3")
Style-AsDate $wsZh.Range("K3")

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("A3").Value = $newFile
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "False"
$wsDe.Range("G3").Value = "b9ee350c-e115-4378-a1b2-248be2ee1af7.815491ef7fccedc18cdde3a25d641bd4f82b7eda.de-de.xlf"
$wsDe.Range("H3").Value = "2016-08-16 08:39:53"
$wsDe.Range("I3").Value = ""
$wsDe.Range("J3").Value = ""
$wsDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDe.Range("L3").Value = ""
$wsDe.Range("M3").Value = "True"
$wsDe.Range("N3").Value = ""
$wsDe.Range("O3").Value = "False"
$wsDe.Range("P3").Value = ""

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $newFileUrl, "", "", $newFile) | Out-Null
Style-AsHyperlink $wsDe.Range("A3")
Style-AsDate $wsDe.Range("H3")
Style-AsDate $wsDe.Range("K3")
